$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '35.278.39'
Set-TextValue $ws 'E2' '  -0.73%  '
Set-TextValue $ws 'D3' '1.896.48'
Set-TextValue $ws 'E3' '  -0.85%  '
Set-TextValue $ws 'E4' '  -0.07%  '
Set-TextValue $ws 'D5' '246.35'
Set-TextValue $ws 'E5' '  -0.19%  '
Set-TextValue $ws 'E6' '  +9.23%  '
Set-TextValue $ws 'E7' '  -0.04%  '
Set-TextValue $ws 'D8' '40.43'
Set-TextValue $ws 'E8' '  -4.06%  '
Set-TextValue $ws 'D9' '0.348'
Set-TextValue $ws 'E9' '  +2.55%  '
Set-TextValue $ws 'D10' '51.92'
Set-TextValue $ws 'E10' '  +6.54%  '
Set-TextValue $ws 'D11' '0.0722'
Set-TextValue $ws 'E11' '  +2.20%  '
Set-TextValue $ws 'D12' '0.0988'
Set-TextValue $ws 'E12' '  -1.00%  '
Set-TextValue $ws 'D13' '2.171.28'
Set-TextValue $ws 'E13' '  -0.71%  '
Set-TextValue $ws 'E14' '  +0.63%  '
Set-TextValue $ws 'E15' '  +2.18%  '
Set-TextValue $ws 'D16' '1.890.01'
Set-TextValue $ws 'E16' '  -1.03%  '
Set-TextValue $ws 'D17' '4.83'
Set-TextValue $ws 'E17' '  -1.06%  '
Set-TextValue $ws 'D18' '35.254.58'
Set-TextValue $ws 'E18' '  -0.85%  '
Set-TextValue $ws 'D19' '72.53'
Set-TextValue $ws 'E19' '  +0.54%  '
Set-TextValue $ws 'D20' '0.0₃0819'
Set-TextValue $ws 'E20' '  -0.20%  '
Set-TextValue $ws 'D21' '240.99'
Set-TextValue $ws 'E21' '  -1.18%  '
Set-TextValue $ws 'D22' '12.75'
Set-TextValue $ws 'E22' '  +1.44%  '
Set-TextValue $ws 'D23' '4.87'
Set-TextValue $ws 'E23' '  -0.58%  '
Set-TextValue $ws 'E24' '  +0.00%  '
Set-TextValue $ws 'E25' '  +1.10%  '
Set-TextValue $ws 'D26' '2.34'
Set-TextValue $ws 'E26' '  +6.19%  '
Set-TextValue $ws 'D27' '167.73'
Set-TextValue $ws 'E27' '  -2.74%  '
Set-TextValue $ws 'D28' '8.56'
Set-TextValue $ws 'E28' '  -0.20%  '
Set-TextValue $ws 'E29' '  +5.36%  '
Set-TextValue $ws 'E30' '  +3.83%  '
Set-TextValue $ws 'D32' '4.18'
Set-TextValue $ws 'E32' '  +1.51%  '
Set-TextValue $ws 'D33' '0.0570'
Set-TextValue $ws 'E33' '  -0.33%  '
Set-TextValue $ws 'E34' '  -0.11%  '
Set-TextValue $ws 'D35' '1.86'
Set-TextValue $ws 'E35' '  +6.69%  '
Set-TextValue $ws 'E36' '  -1.91%  '
Set-TextValue $ws 'D37' '0.909'
Set-TextValue $ws 'E37' '  -6.36%  '
Set-TextValue $ws 'D38' '1.48'
Set-TextValue $ws 'E38' '  +5.86%  '
Set-TextValue $ws 'E39' '  -0.21%  '
Set-TextValue $ws 'B40' 'Aave'
Set-TextValue $ws 'C40' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D40' '95.44'
Set-TextValue $ws 'E40' '  +4.48%  '
Set-TextValue $ws 'B41' 'Kaspa'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D41' '0.0659'
Set-TextValue $ws 'E41' '  +9.30%  '
Set-TextValue $ws 'B42' 'ARBITRUM'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws 'D42' '1.09'
Set-TextValue $ws 'E42' '  -2.01%  '
Set-TextValue $ws 'B43' 'InjectiveProtocol'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D43' '16.39'
Set-TextValue $ws 'E43' '  +4.27%  '
Set-TextValue $ws 'D44' '0.0208'
Set-TextValue $ws 'E44' '  +0.91%  '
Set-TextValue $ws 'D45' '1.353.83'
Set-TextValue $ws 'E45' '  -0.13%  '
Set-TextValue $ws 'E46' '  +0.43%  '
Set-TextValue $ws 'E47' '  +0.03%  '
Set-TextValue $ws 'B48' 'MXToken'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D48' '2.79'
Set-TextValue $ws 'E48' '  +0.59%  '
Set-TextValue $ws 'B49' 'Gas'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
Set-TextValue $ws 'D49' '12.42'
Set-TextValue $ws 'E49' '  -1.60%  '
Set-TextValue $ws 'D50' '45.45'
Set-TextValue $ws 'E50' '  -8.14%  '
Set-TextValue $ws 'E51' '  -2.91%  '

Write-Output "done"